# Replaced SparkFun Electret microphone with SparkFun MEMS microphone.
# Added notes and clarifications to parts list.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 29: P10 COZIR CO2 sensor part number + new note
$ws.Range("F29").Value = "CozIR-A"
$ws.Range("M29").Value = "Used 5,000ppm version (w/o temperature & humidity sensor), but any of the CoZIR-A sensors work."

# Row 26: P1 microphone designator/designation/part number/Pololu part number
$ws.Range("F26").Value = "BOB-09868"
$ws.Range("B26").Value = "INMP401"
$ws.Range("D26").Value = "INMP401"
$ws.Range("K26").Value = 1618
